$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The results table was regenerated from updated AutoML run data. The source
# the numbers were pulled from had its "±" character mangled into the two
# character sequence "Â±" (a classic UTF-8 double-encoding bug), so every
# "f1_score_weighted" / "training_time" / "test_time" cell (not just the
# brand new ones) now shows "Â±" instead of "±". "missing_runs" (always the
# literal "[]") and the already-correct "best_seed" numbers are unaffected.

# --- row 2 (4intelligence) ---------------------------------------------------
$ws.Range("B2").Value = "0.258 (0.238 Â± 0.019)"
$ws.Range("C2").Value = "00:06:10 (00:27:23 Â± 00:13:09)"
$ws.Range("D2").Value = "00:00:00 (00:00:02 Â± 00:00:01)"

# --- row 3 (autogluon) — brand new result row --------------------------------
# A3 already holds "autogluon" (unchanged).
$ws.Range("B3").Value = "0.229 (0.191 Â± 0.021)"
$ws.Range("C3").Value = "00:05:03 (00:06:33 Â± 00:01:07)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("E3").Value = "[]"

# F3 ("best_seed") is stored as text ("59"), like the other best_seed cells
# in this column (F4, F5, ...), not as a number. Writing it with a leading
# apostrophe forces text entry; resetting the style back to "Normal"
# afterwards drops the quote-prefix formatting flag again so the cell keeps
# the sheet's default (unstyled) look, matching the other best_seed cells.
$ws.Range("F3").Value = "'59"
$ws.Range("F3").Style = "Normal"

# --- row 4 (autokeras) --------------------------------------------------------
$ws.Range("B4").Value = "0.209 (0.173 Â± 0.018)"
$ws.Range("C4").Value = "00:01:00 (00:01:21 Â± 00:00:19)"
$ws.Range("D4").Value = "00:00:00 (00:00:00 Â± 00:00:00)"

# --- row 5 (autopytorch) ------------------------------------------------------
$ws.Range("B5").Value = "0.217 (0.171 Â± 0.023)"
$ws.Range("C5").Value = "00:05:08 (00:05:16 Â± 00:00:04)"
$ws.Range("D5").Value = "00:00:01 (00:00:02 Â± 00:00:01)"

# --- row 6 (autosklearn) ------------------------------------------------------
$ws.Range("B6").Value = "0.230 (0.198 Â± 0.018)"
$ws.Range("C6").Value = "00:04:57 (00:05:00 Â± 00:00:03)"
$ws.Range("D6").Value = "00:00:01 (00:00:03 Â± 00:00:02)"

# --- row 8 (flaml) -------------------------------------------------------------
$ws.Range("B8").Value = "0.224 (0.157 Â± 0.050)"
$ws.Range("C8").Value = "00:05:00 (00:05:07 Â± 00:00:15)"
$ws.Range("D8").Value = "00:00:00 (00:00:00 Â± 00:00:00)"

# --- row 10 (h2o) ---------------------------------------------------------------
$ws.Range("B10").Value = "0.180 (0.106 Â± 0.044)"
$ws.Range("C10").Value = "00:05:16 (00:05:54 Â± 00:00:26)"
$ws.Range("D10").Value = "00:00:00 (00:00:00 Â± 00:00:00)"

# --- row 12 (pycaret) ------------------------------------------------------------
$ws.Range("B12").Value = "0.224 (0.191 Â± 0.017)"
$ws.Range("C12").Value = "00:09:28 (00:09:48 Â± 00:00:11)"
$ws.Range("D12").Value = "00:00:00 (00:00:00 Â± 00:00:00)"

# --- row 13 (tpot) -----------------------------------------------------------------
$ws.Range("B13").Value = "0.213 (0.184 Â± 0.017)"
$ws.Range("C13").Value = "00:05:18 (00:06:48 Â± 00:01:27)"
$ws.Range("D13").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
